# edit.ps1 - Applies the "Modifications fichiers semaine 3" commit to CDC.docx
# using Word COM-interop automation (iron_native runtime).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits right after "Responsable"
#    near the top of the document (it is relocated later in the document,
#    see step 8 below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. "          TCP" -> "          OPC"  (bold heading line)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TCP", $true, $true, $false, $false, $false, $true, 1, $false, "OPC", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. ". Les documentations concernant leurs travaux (spéculations) "
#    -> ". Les documentations concernant ses travaux (spécifications générales) "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("leurs travaux", $true, $false, $false, $false, $false, $true, 1, $false, "ses travaux", 2) | Out-Null
$d.Content.Find.Execute("(spéculations)", $true, $false, $false, $false, $false, $true, 1, $false, "(spécifications générales)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Add "M. Réau Pierrick (service informatique)," before "M. Guilloteau Kevin"
#    ", M. Guilloteau Kevin" -> ", M. Réau Pierrick (service informatique), M. Guilloteau Kevin"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(", M. Guilloteau Kevin", $true, $false, $false, $false, $false, $true, 1, $false, ", M. Réau Pierrick (service informatique), M. Guilloteau Kevin", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "modifie" -> "modifier"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("modifie", $true, $true, $false, $false, $false, $true, 1, $false, "modifier", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. "nous voulons que" -> "nous souhaitons que", then re-insert the "_GoBack"
#    bookmark right after "souhaitons" (its new location in the document).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("voulons", $true, $true, $false, $false, $false, $true, 1, $false, "souhaitons", 2) | Out-Null

$rng = $d.Content
$found = $rng.Find.Execute("souhaitons", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)
}

# ---------------------------------------------------------------------------
# 7. "pour/avant" -> "pour"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("pour/avant", $true, $false, $false, $false, $false, $true, 1, $false, "pour", 2) | Out-Null
